# Replace the original single paragraph (bold "Costo:" / DOCPROPERTY fields /
# "n:" / "Importo da fatturare:" laid out with tabs) with a 5-column,
# 1-row borderless table (style "Grigliatabella" / Table Grid) that holds
# the same DOCPROPERTY fields in dedicated, right-justified cells, plus a
# bold "N:" label cell and a bold "Totale" label cell. A trailing empty
# paragraph is left after the table, matching the target layout.

$d = $word.ActiveDocument

# The whole (and only) paragraph in the document is what gets replaced by
# the table. Grab its Range before mutating anything.
$rng = $d.Paragraphs(1).Range

# Body fragment for the new table + following blank paragraph. Word
# consumes the very last paragraph mark of the inserted fragment to keep
# the document's own trailing/sectPr-holding paragraph alive, so two
# trailing <w:p/> are supplied to end up with exactly one visible blank
# paragraph after the table.
$body = @'
<w:tbl><w:tblPr><w:tblStyle w:val="Grigliatabella"/><w:tblW w:w="8330" w:type="dxa"/><w:tblBorders><w:top w:val="none" w:sz="0" w:space="0" w:color="auto"/><w:left w:val="none" w:sz="0" w:space="0" w:color="auto"/><w:bottom w:val="none" w:sz="0" w:space="0" w:color="auto"/><w:right w:val="none" w:sz="0" w:space="0" w:color="auto"/><w:insideH w:val="none" w:sz="0" w:space="0" w:color="auto"/><w:insideV w:val="none" w:sz="0" w:space="0" w:color="auto"/></w:tblBorders><w:tblLayout w:type="fixed"/><w:tblLook w:val="04A0"/></w:tblPr><w:tblGrid><w:gridCol w:w="2518"/><w:gridCol w:w="1134"/><w:gridCol w:w="1559"/><w:gridCol w:w="1560"/><w:gridCol w:w="1559"/></w:tblGrid><w:tr><w:tc><w:tcPr><w:tcW w:w="2518" w:type="dxa"/></w:tcPr><w:p><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:fldChar w:fldCharType="begin"/></w:r><w:r><w:instrText xml:space="preserve"> DOCPROPERTY  Cost.Description  \* MERGEFORMAT </w:instrText></w:r><w:r><w:fldChar w:fldCharType="separate"/></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Cost.Description</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:fldChar w:fldCharType="end"/></w:r><w:r><w:t xml:space="preserve">  </w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="1134" w:type="dxa"/></w:tcPr><w:p><w:pPr><w:jc w:val="right"/><w:rPr><w:b/></w:rPr></w:pPr><w:r><w:rPr><w:b/></w:rPr><w:t>N:</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="1559" w:type="dxa"/></w:tcPr><w:p><w:pPr><w:jc w:val="right"/></w:pPr><w:r><w:fldChar w:fldCharType="begin"/></w:r><w:r><w:instrText xml:space="preserve"> DOCPROPERTY  Cost.Quantity  \* MERGEFORMAT </w:instrText></w:r><w:r><w:fldChar w:fldCharType="separate"/></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Cost.Quantity</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:fldChar w:fldCharType="end"/></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="1560" w:type="dxa"/></w:tcPr><w:p><w:pPr><w:jc w:val="right"/><w:rPr><w:b/></w:rPr></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:b/></w:rPr><w:t>Totale</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="1559" w:type="dxa"/></w:tcPr><w:p><w:pPr><w:jc w:val="right"/></w:pPr><w:r><w:fldChar w:fldCharType="begin"/></w:r><w:r><w:instrText xml:space="preserve"> DOCPROPERTY  Cost.GranTotalCost  \* MERGEFORMAT </w:instrText></w:r><w:r><w:fldChar w:fldCharType="separate"/></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Cost.GranTotalCost</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:fldChar w:fldCharType="end"/></w:r></w:p></w:tc></w:tr></w:tbl><w:p/><w:p/>
'@

$xml = '<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' + $body + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

[void]$rng.InsertXML($xml)

Write-Host "Table inserted; tables=$($d.Tables.Count) paragraphs=$($d.Paragraphs.Count)"
